$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 272-293 (D, L, M, N, O, P, Q, R, S, T) ---
# Row 272
$ws.Cells.Item(272, 4).Value = 45132
$ws.Cells.Item(272, 12).Value = "Extra (doble especial)"
$ws.Cells.Item(272, 13).Value = 370
$ws.Cells.Item(272, 14).Value = 21600
$ws.Cells.Item(272, 15).Value = 24000
$ws.Cells.Item(272, 16).Value = 23027
$ws.Cells.Item(272, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(272, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(272, 19).Value = 2878
$ws.Cells.Item(272, 20).Value = 8

# Row 273
$ws.Cells.Item(273, 4).Value = 45132
$ws.Cells.Item(273, 12).Value = "Primera"
$ws.Cells.Item(273, 13).Value = 200
$ws.Cells.Item(273, 14).Value = 17600
$ws.Cells.Item(273, 15).Value = 17600
$ws.Cells.Item(273, 16).Value = 17600
$ws.Cells.Item(273, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(273, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(273, 19).Value = 2200
$ws.Cells.Item(273, 20).Value = 8

# Row 274
$ws.Cells.Item(274, 4).Value = 45132
$ws.Cells.Item(274, 12).Value = "Segunda"
$ws.Cells.Item(274, 13).Value = 180
$ws.Cells.Item(274, 14).Value = 14400
$ws.Cells.Item(274, 15).Value = 14400
$ws.Cells.Item(274, 16).Value = 14400
$ws.Cells.Item(274, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(274, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(274, 19).Value = 1800
$ws.Cells.Item(274, 20).Value = 8

# Row 275
$ws.Cells.Item(275, 4).Value = 44495
$ws.Cells.Item(275, 12).Value = "Cuarta"
$ws.Cells.Item(275, 13).Value = 320
$ws.Cells.Item(275, 14).Value = 1200
$ws.Cells.Item(275, 15).Value = 1200
$ws.Cells.Item(275, 16).Value = 1200
$ws.Cells.Item(275, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(275, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(275, 19).Value = 1200
$ws.Cells.Item(275, 20).Value = 1

# Row 276
$ws.Cells.Item(276, 4).Value = 44495
$ws.Cells.Item(276, 12).Value = "Especial"
$ws.Cells.Item(276, 13).Value = 280
$ws.Cells.Item(276, 14).Value = 20000
$ws.Cells.Item(276, 15).Value = 20000
$ws.Cells.Item(276, 16).Value = 20000
$ws.Cells.Item(276, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(276, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(276, 19).Value = 2500
$ws.Cells.Item(276, 20).Value = 8

# Row 277
$ws.Cells.Item(277, 4).Value = 44495
$ws.Cells.Item(277, 12).Value = "Extra (doble especial)"
$ws.Cells.Item(277, 13).Value = 300
$ws.Cells.Item(277, 14).Value = 24000
$ws.Cells.Item(277, 15).Value = 24000
$ws.Cells.Item(277, 16).Value = 24000
$ws.Cells.Item(277, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(277, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(277, 19).Value = 3000
$ws.Cells.Item(277, 20).Value = 8

# Row 278
$ws.Cells.Item(278, 4).Value = 44495
$ws.Cells.Item(278, 12).Value = "Primera"
$ws.Cells.Item(278, 13).Value = 350
$ws.Cells.Item(278, 14).Value = 16000
$ws.Cells.Item(278, 15).Value = 16000
$ws.Cells.Item(278, 16).Value = 16000
$ws.Cells.Item(278, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(278, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(278, 19).Value = 2000
$ws.Cells.Item(278, 20).Value = 8

# Row 279
$ws.Cells.Item(279, 4).Value = 44495
$ws.Cells.Item(279, 12).Value = "Segunda"
$ws.Cells.Item(279, 13).Value = 310
$ws.Cells.Item(279, 14).Value = 14400
$ws.Cells.Item(279, 15).Value = 14400
$ws.Cells.Item(279, 16).Value = 14400
$ws.Cells.Item(279, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(279, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(279, 19).Value = 1800
$ws.Cells.Item(279, 20).Value = 8

# Row 280
$ws.Cells.Item(280, 4).Value = 44495
$ws.Cells.Item(280, 12).Value = "Tercera"
$ws.Cells.Item(280, 13).Value = 250
$ws.Cells.Item(280, 14).Value = 1500
$ws.Cells.Item(280, 15).Value = 1500
$ws.Cells.Item(280, 16).Value = 1500
$ws.Cells.Item(280, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(280, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(280, 19).Value = 1500
$ws.Cells.Item(280, 20).Value = 1

# Row 281
$ws.Cells.Item(281, 4).Value = 44511
$ws.Cells.Item(281, 12).Value = "Cuarta"
$ws.Cells.Item(281, 13).Value = 250
$ws.Cells.Item(281, 14).Value = 1200
$ws.Cells.Item(281, 15).Value = 1200
$ws.Cells.Item(281, 16).Value = 1200
$ws.Cells.Item(281, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(281, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(281, 19).Value = 1200
$ws.Cells.Item(281, 20).Value = 1

# Row 282
$ws.Cells.Item(282, 4).Value = 44511
$ws.Cells.Item(282, 12).Value = "Especial"
$ws.Cells.Item(282, 13).Value = 280
$ws.Cells.Item(282, 14).Value = 20000
$ws.Cells.Item(282, 15).Value = 20000
$ws.Cells.Item(282, 16).Value = 20000
$ws.Cells.Item(282, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(282, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(282, 19).Value = 2500
$ws.Cells.Item(282, 20).Value = 8

# Row 283
$ws.Cells.Item(283, 4).Value = 44511
$ws.Cells.Item(283, 12).Value = "Extra (doble especial)"
$ws.Cells.Item(283, 13).Value = 350
$ws.Cells.Item(283, 14).Value = 24000
$ws.Cells.Item(283, 15).Value = 24000
$ws.Cells.Item(283, 16).Value = 24000
$ws.Cells.Item(283, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(283, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(283, 19).Value = 3000
$ws.Cells.Item(283, 20).Value = 8

# Row 284
$ws.Cells.Item(284, 4).Value = 44511
$ws.Cells.Item(284, 12).Value = "Primera"
$ws.Cells.Item(284, 13).Value = 330
$ws.Cells.Item(284, 14).Value = 16000
$ws.Cells.Item(284, 15).Value = 16000
$ws.Cells.Item(284, 16).Value = 16000
$ws.Cells.Item(284, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(284, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(284, 19).Value = 2000
$ws.Cells.Item(284, 20).Value = 8

# Row 285
$ws.Cells.Item(285, 4).Value = 44511
$ws.Cells.Item(285, 12).Value = "Segunda"
$ws.Cells.Item(285, 13).Value = 300
$ws.Cells.Item(285, 14).Value = 14400
$ws.Cells.Item(285, 15).Value = 14400
$ws.Cells.Item(285, 16).Value = 14400
$ws.Cells.Item(285, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(285, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(285, 19).Value = 1800
$ws.Cells.Item(285, 20).Value = 8

# Row 286
$ws.Cells.Item(286, 4).Value = 44391
$ws.Cells.Item(286, 12).Value = "Tercera"
$ws.Cells.Item(286, 13).Value = 220
$ws.Cells.Item(286, 14).Value = 1400
$ws.Cells.Item(286, 15).Value = 1400
$ws.Cells.Item(286, 16).Value = 1400
$ws.Cells.Item(286, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(286, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(286, 19).Value = 1400
$ws.Cells.Item(286, 20).Value = 1

# Row 287
$ws.Cells.Item(287, 4).Value = 44391
$ws.Cells.Item(287, 12).Value = "Especial"
$ws.Cells.Item(287, 13).Value = 6
$ws.Cells.Item(287, 14).Value = 2500
$ws.Cells.Item(287, 15).Value = 2500
$ws.Cells.Item(287, 16).Value = 2500
$ws.Cells.Item(287, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(287, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(287, 19).Value = 2500
$ws.Cells.Item(287, 20).Value = 1

# Row 288
$ws.Cells.Item(288, 4).Value = 44391
$ws.Cells.Item(288, 12).Value = "Extra (doble especial)"
$ws.Cells.Item(288, 13).Value = 3
$ws.Cells.Item(288, 14).Value = 3000
$ws.Cells.Item(288, 15).Value = 3000
$ws.Cells.Item(288, 16).Value = 3000
$ws.Cells.Item(288, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(288, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(288, 19).Value = 3000
$ws.Cells.Item(288, 20).Value = 1

# Row 289
$ws.Cells.Item(289, 4).Value = 44391
$ws.Cells.Item(289, 12).Value = "Primera"
$ws.Cells.Item(289, 13).Value = 7
$ws.Cells.Item(289, 14).Value = 2000
$ws.Cells.Item(289, 15).Value = 2000
$ws.Cells.Item(289, 16).Value = 2000
$ws.Cells.Item(289, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(289, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(289, 19).Value = 2000
$ws.Cells.Item(289, 20).Value = 1

# Row 290
$ws.Cells.Item(290, 4).Value = 44391
$ws.Cells.Item(290, 12).Value = "Segunda"
$ws.Cells.Item(290, 13).Value = 10
$ws.Cells.Item(290, 14).Value = 1500
$ws.Cells.Item(290, 15).Value = 1500
$ws.Cells.Item(290, 16).Value = 1500
$ws.Cells.Item(290, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(290, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(290, 19).Value = 1500
$ws.Cells.Item(290, 20).Value = 1

# Row 291
$ws.Cells.Item(291, 4).Value = 44859
$ws.Cells.Item(291, 12).Value = "Cuarta"
$ws.Cells.Item(291, 13).Value = 180
$ws.Cells.Item(291, 14).Value = 1500
$ws.Cells.Item(291, 15).Value = 1500
$ws.Cells.Item(291, 16).Value = 1500
$ws.Cells.Item(291, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(291, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(291, 19).Value = 1500
$ws.Cells.Item(291, 20).Value = 1

# Row 292
$ws.Cells.Item(292, 4).Value = 44859
$ws.Cells.Item(292, 12).Value = "Especial"
$ws.Cells.Item(292, 13).Value = 250
$ws.Cells.Item(292, 14).Value = 22400
$ws.Cells.Item(292, 15).Value = 22400
$ws.Cells.Item(292, 16).Value = 22400
$ws.Cells.Item(292, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(292, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(292, 19).Value = 2800
$ws.Cells.Item(292, 20).Value = 8

# Row 293
$ws.Cells.Item(293, 4).Value = 44859
$ws.Cells.Item(293, 12).Value = "Extra (doble especial)"
$ws.Cells.Item(293, 13).Value = 260
$ws.Cells.Item(293, 14).Value = 24000
$ws.Cells.Item(293, 15).Value = 24000
$ws.Cells.Item(293, 16).Value = 24000
$ws.Cells.Item(293, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(293, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(293, 19).Value = 3000
$ws.Cells.Item(293, 20).Value = 8

# --- Add new rows 294-296 ---
# Row 294
$ws.Cells.Item(294, 1).Value = 9
$ws.Cells.Item(294, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(294, 3).Value = "Metropolitana"
$ws.Cells.Item(294, 4).Value = 44859
$ws.Cells.Item(294, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(294, 5).Value = 13
$ws.Cells.Item(294, 6).Value = "Fruta"
$ws.Cells.Item(294, 7).Value = 100107
$ws.Cells.Item(294, 8).Value = "Otros"
$ws.Cells.Item(294, 9).Value = 100107002
$ws.Cells.Item(294, 10).Value = "Chirimoya"
$ws.Cells.Item(294, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(294, 12).Value = "Primera"
$ws.Cells.Item(294, 13).Value = 220
$ws.Cells.Item(294, 14).Value = 20800
$ws.Cells.Item(294, 15).Value = 20800
$ws.Cells.Item(294, 16).Value = 20800
$ws.Cells.Item(294, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(294, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(294, 19).Value = 2600
$ws.Cells.Item(294, 20).Value = 8

# Row 295
$ws.Cells.Item(295, 1).Value = 9
$ws.Cells.Item(295, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(295, 3).Value = "Metropolitana"
$ws.Cells.Item(295, 4).Value = 44859
$ws.Cells.Item(295, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(295, 5).Value = 13
$ws.Cells.Item(295, 6).Value = "Fruta"
$ws.Cells.Item(295, 7).Value = 100107
$ws.Cells.Item(295, 8).Value = "Otros"
$ws.Cells.Item(295, 9).Value = 100107002
$ws.Cells.Item(295, 10).Value = "Chirimoya"
$ws.Cells.Item(295, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(295, 12).Value = "Segunda"
$ws.Cells.Item(295, 13).Value = 200
$ws.Cells.Item(295, 14).Value = 17600
$ws.Cells.Item(295, 15).Value = 17600
$ws.Cells.Item(295, 16).Value = 17600
$ws.Cells.Item(295, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(295, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(295, 19).Value = 2200
$ws.Cells.Item(295, 20).Value = 8

# Row 296
$ws.Cells.Item(296, 1).Value = 9
$ws.Cells.Item(296, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(296, 3).Value = "Metropolitana"
$ws.Cells.Item(296, 4).Value = 44859
$ws.Cells.Item(296, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(296, 5).Value = 13
$ws.Cells.Item(296, 6).Value = "Fruta"
$ws.Cells.Item(296, 7).Value = 100107
$ws.Cells.Item(296, 8).Value = "Otros"
$ws.Cells.Item(296, 9).Value = 100107002
$ws.Cells.Item(296, 10).Value = "Chirimoya"
$ws.Cells.Item(296, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(296, 12).Value = "Tercera"
$ws.Cells.Item(296, 13).Value = 200
$ws.Cells.Item(296, 14).Value = 1800
$ws.Cells.Item(296, 15).Value = 1800
$ws.Cells.Item(296, 16).Value = 1800
$ws.Cells.Item(296, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(296, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(296, 19).Value = 1800
$ws.Cells.Item(296, 20).Value = 1

# --- Update sheet dimension reference ---
Write-Host "Done. UsedRange:" $ws.UsedRange.Address()